$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1
$ws.Range("D1").Value = "Oct. Days Deployed per Week"
$ws.Range("E1").Value = "Sep. Days Deployed per Week"

# Move the active selection to E1 (was E9)
$ws.Activate()
$ws.Range("E1").Select()
